# A new weekly price record was inserted as the new row 432 ("Segunda"
# quality, week of 2022-07-12 / serial 44754), pushing every existing
# record from row 432 down through row 533 one row further (to rows
# 433-534). This mirrors inserting one row above the prior row 432 and
# filling it with the new data, which is exactly what Excel does when a
# new daily/weekly observation is spliced into the middle of this
# chronological price log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 432; everything from the old row 432
# onward (through the old last row, 533) shifts down by one row.
$ws.Rows.Item(432).Insert()

# Populate the newly inserted row 432 with the new record's data.
$ws.Cells.Item(432, 1).Value  = 6
$ws.Cells.Item(432, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(432, 3).Value  = "Metropolitana"
$ws.Cells.Item(432, 4).Value  = 44754
$ws.Cells.Item(432, 5).Value  = 13
$ws.Cells.Item(432, 6).Value  = 100112044
$ws.Cells.Item(432, 7).Value  = "Perejil"
$ws.Cells.Item(432, 8).Value  = "Sin especificar"
$ws.Cells.Item(432, 9).Value  = "Segunda"
$ws.Cells.Item(432, 10).Value = 70
$ws.Cells.Item(432, 11).Value = 19000
$ws.Cells.Item(432, 12).Value = 20000
$ws.Cells.Item(432, 13).Value = 19571
$ws.Cells.Item(432, 14).Value = "`$/docena de atados"
$ws.Cells.Item(432, 15).Value = "Región Metropolitana"
$ws.Cells.Item(432, 16).Value = 6524
$ws.Cells.Item(432, 17).Value = 3
$ws.Cells.Item(432, 18).Value = "Hortaliza"
